$wb = $excel.ActiveWorkbook

$belgium = $wb.Worksheets.Item("Belgium")

# Create Denmark, Sweden, Norway sheets after Belgium by copying Belgium as a template
$belgium.Copy([System.Reflection.Missing]::Value, $belgium)
$denmark = $wb.Worksheets.Item($belgium.Index + 1)
$denmark.Name = "Denmark"

$denmark.Copy([System.Reflection.Missing]::Value, $denmark)
$sweden = $wb.Worksheets.Item($denmark.Index + 1)
$sweden.Name = "Sweden"

$sweden.Copy([System.Reflection.Missing]::Value, $sweden)
$norway = $wb.Worksheets.Item($sweden.Index + 1)
$norway.Name = "Norway"

# Populate Denmark
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2003"
$denmark.Range("B4").ClearFormats() | Out-Null

# Populate Sweden
$sweden.Range("B2").Value = "Sweden market"
$sweden.Range("B4").Value = "NGC-3465/T2029"

# Populate Norway
$norway.Range("B2").Value = "Norway market"
$norway.Range("B4").Value = "NGC-3464/T1918"

# Selections: Belgium, Denmark, Sweden get "select all" (Ctrl+A); Norway gets B6 as active selection
$belgium.Activate()
$belgium.Cells.Select() | Out-Null

$denmark.Activate()
$denmark.Cells.Select() | Out-Null

$sweden.Activate()
$sweden.Cells.Select() | Out-Null

$norway.Activate()
$norway.Range("B6").Select() | Out-Null
